# Trade #97 closed at 2026-02-17 15:57:26 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.91   # Current Capital
$summary.Range("B4").Value = -0.1      # Total P&L $
$summary.Range("B6").Value = 97        # Total Trades
$summary.Range("B8").Value = 49        # Losing Trades
$summary.Range("B9").Value = 37.11     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.91      # Capital
$status.Range("D4").Value = 97         # Trades
$status.Range("E4").Value = -0.1       # P&L $
$status.Range("F4").Value = -0.09      # P&L %
$status.Range("G4").Value = 37.11      # Win Rate %

# ---------------------------------------------------------------------
# New closed trade (#97) appended to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------
$newRow = @{
    A = 97
    B = "2026-02-17"
    C = "15:57:19"
    D = "MarketMaking"
    E = "UP"
    F = 0.17
    G = 0.16
    H = "CLOSED"
    I = -5.8824
    J = -0.01
    K = 99.91
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 98

    $ws.Cells.Item($row, 1).Value = $newRow.A

    # Force the Date column to stay plain text ("2026-02-17") instead of
    # being auto-converted to a date serial number by Excel's input
    # parser; reset back to the Normal style afterwards so no stray
    # number-format is left applied to the cell.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $newRow.B
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = $newRow.C
    $ws.Cells.Item($row, 4).Value = $newRow.D
    $ws.Cells.Item($row, 5).Value = $newRow.E
    $ws.Cells.Item($row, 6).Value = $newRow.F
    $ws.Cells.Item($row, 7).Value = $newRow.G
    $ws.Cells.Item($row, 8).Value = $newRow.H
    $ws.Cells.Item($row, 9).Value = $newRow.I
    $ws.Cells.Item($row, 10).Value = $newRow.J
    $ws.Cells.Item($row, 11).Value = $newRow.K
    $ws.Cells.Item($row, 12).Value = $newRow.L
    $ws.Cells.Item($row, 13).Value = $newRow.M
    $ws.Cells.Item($row, 14).Value = $newRow.N
    $ws.Cells.Item($row, 15).Value = $newRow.O
    $ws.Cells.Item($row, 16).Value = $newRow.P
    $ws.Cells.Item($row, 17).Value = $newRow.Q
}
